$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:123 down to 44:124
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record's data
$ws.Cells.Item(43, 1).Value = 4
$ws.Cells.Item(43, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(43, 3).Value = "Los Lagos"
$ws.Cells.Item(43, 4).Value = 44757
$ws.Cells.Item(43, 5).Value = 10
$ws.Cells.Item(43, 6).Value = 100112022
$ws.Cells.Item(43, 7).Value = "Arveja Verde"
$ws.Cells.Item(43, 8).Value = "Perfection"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 70
$ws.Cells.Item(43, 11).Value = 44000
$ws.Cells.Item(43, 12).Value = 44000
$ws.Cells.Item(43, 13).Value = 44000
$ws.Cells.Item(43, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(43, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(43, 16).Value = 1760
$ws.Cells.Item(43, 17).Value = 25
$ws.Cells.Item(43, 18).Value = "Hortaliza"
